# Fruta / hortaliza, semanal
# Update the weekly price rows: the sheet rows for "Bruselas (repollito)" were
# re-fetched/re-shuffled for several dates. Row 3 and row 8 stay unchanged;
# rows 2,4,5,6,7,9,10,11,12,13 get new Fecha/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44398; J = 130; K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 4;  D = 44435; J = 140; K = 21000; L = 23000; M = 21714; P = 1448 },
    @{ Row = 5;  D = 44449; J = 220; K = 22000; L = 24000; M = 23091; P = 1539 },
    @{ Row = 6;  D = 44446; J = 150; K = 22000; L = 24000; M = 22667; P = 1511 },
    @{ Row = 7;  D = 44406; J = 400; K = 20000; L = 22000; M = 20850; P = 1390 },
    @{ Row = 9;  D = 44399; J = 150; K = 22000; L = 22000; M = 22000; P = 1467 },
    @{ Row = 10; D = 44400; J = 130; K = 24000; L = 24000; M = 24000; P = 1600 },
    @{ Row = 11; D = 44392; J = 220; K = 23000; L = 23000; M = 23000; P = 1533 },
    @{ Row = 12; D = 44453; J = 280; K = 20000; L = 22000; M = 21286; P = 1419 },
    @{ Row = 13; D = 44365; J = 580; K = 20000; L = 22000; M = 21103; P = 1407 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("J$r").Value = $u.J
    $ws.Range("K$r").Value = $u.K
    $ws.Range("L$r").Value = $u.L
    $ws.Range("M$r").Value = $u.M
    $ws.Range("P$r").Value = $u.P
}
